# Deploying to gh-pages: add the "2021" data column (R) to the
# "3.в.1-вакцина" sheet, mirroring the formatting of the existing 2020
# column (Q), adjust header row heights, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height tweaks on the title rows -----------------------------
$ws.Rows.Item(1).RowHeight = 41.25
$ws.Rows.Item(2).RowHeight = 15

# --- New column R: copy the number-format/style from column Q, then ---
# --- stamp in the 2021 values (rows with no data stay format-only) ----
$newColumnValues = @{
    3  = $null
    4  = 2021
    5  = $null
    6  = 88.796593100633856
    7  = 86.908583391486388
    8  = 89.680106631122953
    9  = 95.775910364145659
    10 = 96.517042279754136
    11 = 90.311530128242666
    12 = 90.746324915190343
    13 = 90.894107952204379
    14 = 81.065680730752504
    15 = 85.088888888888889
    16 = $null
    17 = 93.37839883628321
    18 = 93.091416608513612
    19 = 94.815061646117954
    20 = 100.53781512605042
    21 = 100.33525796237662
    22 = 93.78989283832054
    23 = 95.401432340746325
    24 = 92.308748798242007
    25 = 89.338842975206617
    26 = 87.955555555555549
    27 = $null
    28 = 89.631204460036727
    29 = 89.204466154919743
    30 = 84.751749416861045
    31 = 96.201680672268907
    32 = 95.567144719687093
    33 = 91.330444457457389
    34 = 91.368262344515642
    35 = 92.345373803964662
    36 = 88.660287081339717
    37 = 84.944444444444443
    38 = $null
}

foreach ($row in ($newColumnValues.Keys | Sort-Object)) {
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")

    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats

    $value = $newColumnValues[$row]
    if ($null -ne $value) {
        $dstCell.Value = $value
    }
}

# --- Move the active selection to R3, as in the authored workbook -----
$ws.Range("R3").Select()
